$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the PURL base URL (test2 -> test) in B1 and C3
$ws.Range("B1").Value = "http://purl.org/test/variables/"
$ws.Range("C3").Value = "http://purl.org/test/variables/"

# Clear placeholder/example metadata values
$ws.Range("B10").Value = ""
$ws.Range("B11").Value = ""
$ws.Range("B12").Value = ""

# Row 19: rename term Test -> test, clear description
$ws.Range("A19").Value = "vars:test"
$ws.Range("B19").Value = "test"
$ws.Range("E19").Value = ""

# Row 20: rename term Computerscientist -> working, clear description and related term
$ws.Range("A20").Value = "vars:working"
$ws.Range("B20").Value = "working"
$ws.Range("E20").Value = ""
$ws.Range("F20").Value = ""

# Row 21: clear term Computerscience, leave only the prefix
$ws.Range("A21").Value = "vars:"
$ws.Range("B21").Value = ""
$ws.Range("E21").Value = ""

$wb.Save()
